$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: replace the old blank B10 cell with a comment line starting in A10 ---
$ws.Range("A10").Value = "Do jeito que está me parece muito simplificado, parece que existem funcionalidades agregadas dentro de outras e/ou funcionalidades faltando, a saber:"
$ws.Range("B10").Clear()

# --- Rows 11-13: bullet points, formatted like the rest of the sheet (Arial 10, no explicit color) ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("B11:B13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("B12").Value = "Login - o texto menciona que o usuário pode alterar sua senha. Não vi esse caso de uso"
$ws.Range("B11").Value = "Cadastro de Usuário - o usuário só pode se cadastrar? Não tem como alterar seus dados depois de cadastrado?"
$ws.Range("B13").Value = "Visualizar itens (um conjunto de itens, normalmente em uma lista), é diferente de visualizar um produto. Necessita ser melhorada essa questão. Normalmente é feita uma consulta, por nome ou categoria, e exibida uma lista de produtos que atendem a pesquisa. Deveria haver uma funcionalidade para isso."

$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 15.75

# Row 13 wraps across the merged B13:F13 block and grows taller
$ws.Range("B13:F13").Merge()
$ws.Range("B13:F13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 39.75

# --- Rows 14-16: further remarks, Arial 10 with explicit black color ---
$ws.Range("B14").Value = "Selecionar Produto é uma coisa e adicionar ao carrinho é outra. Não vi nada para adicionar ao carrinho. Também não vi nada para remover do carrinho"
$ws.Range("B15").Value = "Comprar também é uma funcionalidade muito genérica. E também tem que ter o pagamento como uma funcionalidade separada."
$ws.Range("B16").Value = "Deveria existir também um caso de uso para manter as categorias dos produtos do site"

$rColored = $ws.Range("B14:B16")
$rColored.Font.Name = "Arial"
$rColored.Font.Size = 10
$rColored.Font.Color = 0

$ws.Rows.Item(14).RowHeight = 15.75
$ws.Rows.Item(15).RowHeight = 15.75
$ws.Rows.Item(16).RowHeight = 15.75

# --- Page setup / selection ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("B17").Select() | Out-Null
